$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.745.04"
$ws.Range("E2").Value = "  +3.75%  "

$ws.Range("D3").Value = "2.251.52"
$ws.Range("E3").Value = "  +3.45%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "'253.77"
$ws.Range("E5").Value = "  -0.13%  "

$ws.Range("E6").Value = "  +1.51%  "

$ws.Range("D7").Value = "'71.91"
$ws.Range("E7").Value = "  +6.25%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("D9").Value = "'0.643"
$ws.Range("E9").Value = "  +12.20%  "

$ws.Range("D10").Value = "'41.08"
$ws.Range("E10").Value = "  +9.69%  "

$ws.Range("D11").Value = "'59.54"
$ws.Range("E11").Value = "  +0.96%  "

$ws.Range("D12").Value = "'0.0962"
$ws.Range("E12").Value = "  +3.48%  "

$ws.Range("E13").Value = "  +3.61%  "

$ws.Range("E14").Value = "  +0.50%  "

$ws.Range("D15").Value = "2.590.45"
$ws.Range("E15").Value = "  +3.56%  "

$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.885"
$ws.Range("E16").Value = "  +2.03%  "

$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "'14.80"
$ws.Range("E17").Value = "  +2.80%  "

$ws.Range("D18").Value = "2.251.73"
$ws.Range("E18").Value = "  +3.06%  "

$ws.Range("D19").Value = "42.702.94"
$ws.Range("E19").Value = "  +3.83%  "

$ws.Range("D20").Value = "0.0₃0981"
$ws.Range("E20").Value = "  +3.24%  "

$ws.Range("E21").Value = "  +1.75%  "

$ws.Range("D22").Value = "'72.98"
$ws.Range("E22").Value = "  +1.94%  "

$ws.Range("D23").Value = "'235.42"
$ws.Range("E23").Value = "  +1.79%  "

$ws.Range("E24").Value = "  +5.18%  "

$ws.Range("E25").Value = "  +1.16%  "

$ws.Range("D26").Value = "'11.65"
$ws.Range("E26").Value = "  -0.73%  "

$ws.Range("E27").Value = "  +0.22%  "

$ws.Range("E28").Value = "  -2.72%  "

$ws.Range("D30").Value = "'2.22"
$ws.Range("E30").Value = "  +1.99%  "

$ws.Range("D31").Value = "'168.12"
$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("D32").Value = "'21.03"
$ws.Range("E32").Value = "  +2.18%  "

$ws.Range("D33").Value = "'0.129"
$ws.Range("E33").Value = "  +11.05%  "

$ws.Range("D34").Value = "'6.19"
$ws.Range("E34").Value = "  +13.98%  "

$ws.Range("D35").Value = "'0.0783"
$ws.Range("E35").Value = "  +4.95%  "

$ws.Range("E36").Value = "  +1.53%  "

$ws.Range("D37").Value = "'28.24"
$ws.Range("E37").Value = "  +6.42%  "

$ws.Range("D38").Value = "'4.69"
$ws.Range("E38").Value = "  +1.73%  "

$ws.Range("E39").Value = "  +0.91%  "

$ws.Range("E40").Value = "  +6.25%  "

$ws.Range("E41").Value = "  +4.59%  "

$ws.Range("E42").Value = "  +4.73%  "

$ws.Range("D43").Value = "'12.42"
$ws.Range("E43").Value = "  -0.73%  "

$ws.Range("D44").Value = "'64.60"
$ws.Range("E44").Value = "  +0.55%  "

$ws.Range("D45").Value = "'0.202"
$ws.Range("E45").Value = "  +1.00%  "

$ws.Range("D46").Value = "'4.91"
$ws.Range("E46").Value = "  -2.99%  "

$ws.Range("D47").Value = "'8.91"
$ws.Range("E47").Value = "  +3.27%  "

$ws.Range("E48").Value = "  +1.54%  "

$ws.Range("E49").Value = "  +5.79%  "

$ws.Range("E50").Value = "  -0.13%  "

$ws.Range("D51").Value = "'4.41"
$ws.Range("E51").Value = "  +3.44%  "

# Restore default (no explicit) formatting on cells forced to text via quote-prefix
$ws.Range("D5").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D17").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D51").ClearFormats()
